# Adds the "Planogram" sheet (with its Merchandising/planogram KPI row) to
# the L&T test template, matching the "adding planogram function to logic"
# commit: a new worksheet is appended after "Count", made the active sheet,
# and populated with a KPI header row + one data row.

$wb = $excel.ActiveWorkbook

# --- Add the new sheet after the last existing sheet ("Count") -------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet, 1, $null)
$ws.Name = "Planogram"

# Match the tab colour used by the other KPI sheets (RGB 255,192,0 / FFC000)
$ws.Tab.Color = 49407

# --- Column widths (approximate autosize seen in the authored sheet) -------
$ws.Columns.Item(1).ColumnWidth = 27.93
$ws.Columns.Item(2).ColumnWidth = 49.6
$ws.Columns.Item(3).ColumnWidth = 52.18
$ws.Columns.Item(4).ColumnWidth = 8.86
$ws.Columns.Item(5).ColumnWidth = 8.86
$ws.Columns.Item(6).ColumnWidth = 9.97
$ws.Columns.Item(7).ColumnWidth = 9.72
$ws.Columns.Item(8).ColumnWidth = 9.72
$ws.Columns.Item(9).ColumnWidth = 59.32

# --- Header row (row 1) ------------------------------------------------
$ws.Range("A1").Value = "KPI Name"
$ws.Range("B1").Value = "Atomic KPI Name"
$ws.Range("C1").Value = "Template Display Name"
$ws.Range("D1").Value = "Target"
$ws.Range("E1").Value = "Score "
$ws.Range("F1").Value = "Store Type"
$ws.Range("G1").Value = "Attribute_1"
$ws.Range("H1").Value = "Attribute_2"
$ws.Range("I1").Value = "Comment "

$ws.Range("A1:E1").Interior.Color = 13553360
$ws.Range("F1:H1").Interior.Color = 49407
$ws.Range("I1").Interior.Color = 65535

# --- Data row (row 2) ---------------------------------------------------
$ws.Range("A2").Value = "Merchandising"
$ws.Range("B2").Value = "Cooler Merchandised as per planogram"
$ws.Range("C2").Value = "4 Door Cooler, 3 Door Cooler, 2 Door Cooler, 1 Door Cooler"
$ws.Range("D2").Value = 100
$ws.Range("E2").Value = 6
$ws.Range("F2").Value = "QSR"
$ws.Range("I2").Value = "At least one template should pass the planogram for this KPI to pass"

# --- Page setup: margins + header/footer, matching the other KPI sheets ---
$ws.PageSetup.LeftMargin = 56.7
$ws.PageSetup.RightMargin = 56.7
$ws.PageSetup.TopMargin = 75.8
$ws.PageSetup.BottomMargin = 75.8
$ws.PageSetup.HeaderMargin = 56.7
$ws.PageSetup.FooterMargin = 56.7
$ws.PageSetup.PaperSize = 1
$ws.PageSetup.Orientation = 1
$ws.PageSetup.CenterHeader = '&"Times New Roman,Regular"&12&A'
$ws.PageSetup.CenterFooter = '&"Times New Roman,Regular"&12Page &P'

# --- View: zoom to 58%, select B7, make this the active sheet -------------
$ws.Range("B7").Select()
$excel.ActiveWindow.Zoom = 58

# Record the new selection the user left on the "Price" sheet before
# switching away to the new tab.
$priceSheet = $wb.Worksheets.Item("Price")
$priceSheet.Range("M26").Select()

# Activating "Planogram" last makes it the active/selected tab (activeTab=6)
$ws.Activate()
